# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (switching from Strike# to K, with recalculated std/mean
# and s_vals). Only column G (rows 2-13) changes value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 6
    3  = 5
    4  = 3
    5  = 4
    6  = 4
    7  = 4
    8  = 4
    9  = 6
    10 = 5
    11 = 5
    12 = 1
    13 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
